# chore: update Sheets via scheduled runner
# Re-applies refreshed Universalis price snapshots (currentAveragePrice* / Leve
# profit columns H:N) captured for this workbook's before -> after commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 4 (Leve Item ID 5470)
$ws.Range("H4").Value = 501
$ws.Range("I4").Value = 501
$ws.Range("K4").Value = 501
$ws.Range("M4").Value = -387
# row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 67.5
$ws.Range("I5").Value = 67.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 67.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 47.5
$ws.Range("N5").ClearContents()
# row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
# row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
# row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("M77").ClearContents()
# row 95 (Leve Item ID 18200)
$ws.Range("H95").Value = 38879.6
$ws.Range("J95").Value = 38879.6
$ws.Range("L95").Value = 38879.6
$ws.Range("N95").Value = -44371.6
# row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1541.25
$ws.Range("I132").Value = 1475.7142
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4427.142599999999
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1897.142599999999
$ws.Range("N132").Value = -11060
# row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2564.1428
$ws.Range("I138").Value = 1382.3334
$ws.Range("J138").Value = 3036.8667
$ws.Range("K138").Value = 4147.0002
$ws.Range("L138").Value = 9110.6001
$ws.Range("M138").Value = 992.9997999999996
$ws.Range("N138").Value = -19390.6001

$ws = $wb.Worksheets.Item("ARM")
# row 5 (Leve Item ID 5091)
$ws.Range("H5").Value = 90.833336
$ws.Range("I5").Value = 93.8
$ws.Range("J5").Value = 76
$ws.Range("K5").Value = 93.8
$ws.Range("L5").Value = 76
$ws.Range("M5").Value = 18.2
$ws.Range("N5").Value = -300
# row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 640.3333
$ws.Range("J32").Value = 3599
$ws.Range("L32").Value = 3599
$ws.Range("N32").Value = -4173
# row 135 (Leve Item ID 42016)
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# row 4 (Leve Item ID 5091)
$ws.Range("H4").Value = 90.833336
$ws.Range("I4").Value = 93.8
$ws.Range("J4").Value = 76
$ws.Range("K4").Value = 93.8
$ws.Range("L4").Value = 76
$ws.Range("M4").Value = 21.2
$ws.Range("N4").Value = -306
# row 36 (Leve Item ID 2320)
$ws.Range("H36").Value = 7812
$ws.Range("I36").Value = 7812
$ws.Range("K36").Value = 7812
$ws.Range("M36").Value = -7278

$ws = $wb.Worksheets.Item("CRP")
# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 1551.6666
$ws.Range("I31").Value = 1551.6666
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1551.6666
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1256.6666
$ws.Range("N31").ClearContents()
# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 1551.6666
$ws.Range("I34").Value = 1551.6666
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1551.6666
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1349.6666
$ws.Range("N34").ClearContents()
# row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 13000
$ws.Range("I58").Value = 13000
$ws.Range("J58").Value = 13000
$ws.Range("K58").Value = 13000
$ws.Range("L58").Value = 13000
$ws.Range("M58").Value = -12797
$ws.Range("N58").Value = -13406
# row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 500010000
$ws.Range("I86").Value = 500010000
$ws.Range("K86").Value = 500010000
$ws.Range("M86").Value = -500008877
# row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 500010000
$ws.Range("I89").Value = 500010000
$ws.Range("K89").Value = 2500050000
$ws.Range("M89").Value = -2500044384
# row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 13000
$ws.Range("I136").Value = 13000
$ws.Range("J136").Value = 13000
$ws.Range("K136").Value = 39000
$ws.Range("L136").Value = 39000
$ws.Range("M136").Value = -36450
$ws.Range("N136").Value = -44100

$ws = $wb.Worksheets.Item("CUL")
# row 2 (Leve Item ID 4847)
$ws.Range("H2").Value = 55.5
$ws.Range("J2").Value = 55.5
$ws.Range("L2").Value = 333
$ws.Range("N2").Value = -559
# row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 2536.7646
$ws.Range("I4").Value = 2004.762
$ws.Range("J4").Value = 3396.1538
$ws.Range("K4").Value = 6014.286
$ws.Range("L4").Value = 10188.4614
$ws.Range("M4").Value = -5902.286
$ws.Range("N4").Value = -10412.4614
# row 17 (Leve Item ID 4640)
$ws.Range("H17").Value = 451.4
$ws.Range("J17").Value = 451.4
$ws.Range("L17").Value = 1354.2
$ws.Range("N17").Value = -1692.2
# row 23 (Leve Item ID 4858)
$ws.Range("H23").Value = 130.75
$ws.Range("I23").Value = 332.5
$ws.Range("J23").Value = 63.5
$ws.Range("K23").Value = 997.5
$ws.Range("L23").Value = 190.5
$ws.Range("M23").Value = -762.5
$ws.Range("N23").Value = -660.5
# row 40 (Leve Item ID 4827)
$ws.Range("H40").Value = 485.14285
$ws.Range("I40").Value = 483.5
$ws.Range("J40").Value = 495
$ws.Range("K40").Value = 1934
$ws.Range("L40").Value = 1980
$ws.Range("M40").Value = -1865
$ws.Range("N40").Value = -2118
# row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 2131.25
$ws.Range("J55").Value = 4140
$ws.Range("L55").Value = 12420
$ws.Range("N55").Value = -12774
# row 81 (Leve Item ID 12843)
$ws.Range("H81").Value = 2013
$ws.Range("I81").Value = 2013
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6039
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4916
$ws.Range("N81").ClearContents()
# row 84 (Leve Item ID 12843)
$ws.Range("H84").Value = 2013
$ws.Range("I84").Value = 2013
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 18117
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -12501
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# row 57 (Leve Item ID 2876)
$ws.Range("H57").Value = 33198.8
$ws.Range("J57").Value = 33198.8
$ws.Range("L57").Value = 33198.8
$ws.Range("N57").Value = -34838.8
# row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 250003000
$ws.Range("J70").Value = 250003000
$ws.Range("L70").Value = 250003000
$ws.Range("N70").Value = -250003540
# row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 250003000
$ws.Range("J73").Value = 250003000
$ws.Range("L73").Value = 250003000
$ws.Range("N73").Value = -250004872
# row 134 (Leve Item ID 42064)
$ws.Range("H134").Value = 138333.33
$ws.Range("J134").Value = 138333.33
$ws.Range("L134").Value = 414999.99
$ws.Range("N134").Value = -420069.99

$ws = $wb.Worksheets.Item("LTW")
# row 2 (Leve Item ID 2631)
$ws.Range("H2").Value = 1500
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
# row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 5791.5557
$ws.Range("I16").Value = 3624.75
$ws.Range("J16").Value = 7525
$ws.Range("K16").Value = 3624.75
$ws.Range("L16").Value = 7525
$ws.Range("M16").Value = -3454.75
$ws.Range("N16").Value = -7865

$ws = $wb.Worksheets.Item("WVR")
# row 131 (Leve Item ID 34723)
$ws.Range("H131").Value = 57998
$ws.Range("J131").Value = 57998
$ws.Range("L131").Value = 57998
$ws.Range("N131").Value = -68078
# row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2419.4285
$ws.Range("I132").Value = 2419.4285
$ws.Range("K132").Value = 7258.2855
$ws.Range("M132").Value = -4728.2855
